# Automatische test-sync: 2025-08-06 19:48:50
# Appends a new test-mail log row to the "Logs" sheet and bumps the
# matching category counter on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# --- "Logs" sheet: append row 9 -------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(9, 1).Value = "Bestel je 100 M5-bouten zodra je kan?"
$logs.Cells.Item(9, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(9, 3).Value = "Testmail #1: Bestel je 100 M5-bouten zodra je kan?"
$logs.Cells.Item(9, 4).Value = "Inkoop / Bestellingen"
$logs.Cells.Item(9, 5).Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$logs.Cells.Item(9, 6).Value = "2025-08-06 19:48:14"
$logs.Cells.Item(9, 7).Value = "Ja"
$logs.Cells.Item(9, 8).Value = "Ja"
$logs.Cells.Item(9, 9).Value = "Nee"
$logs.Cells.Item(9, 10).Value = "Nee"

# --- extend the conditional formatting ranges to include the new row -----
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "8")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "9")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- "Dashboard" sheet: bump the "Inkoop / Bestellingen" counter ----------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(3, 2).Value = 3
